$p = $ppt.ActivePresentation

# --- 1. Remove the "About AdventHealth" heading (and the blank line under
#        it) from the "About the hospital" slide (slide 6), shape "TextBox 9".
#        The shape uses AutoFit, so PowerPoint will shrink its height to fit
#        the remaining text automatically.
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(9)
$tr = $sh.TextFrame.TextRange

# Paragraph 1 = "About AdventHealth" (18 chars + paragraph mark)
# Paragraph 2 = "" (0 chars + paragraph mark)
# Together that is the first 20 characters of the text range.
$heading = $tr.Characters(1, 20)
$heading.Delete()

# --- 2. Add a new, blank slide at the end of the deck (position 12) so
#        navigation buttons elsewhere in the deck have somewhere to link to.
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)
